$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 12) of data for year 2021, mirroring the layout of
# the previous row (row 11) so formatting (style, empty string cells for
# columns E and M) is preserved, then overwrite the actual values.
$ws.Range("A11:S11").Copy($ws.Range("A12:S12"))

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 5470.6
$ws.Range("C12").Value = 11139.9
$ws.Range("D12").Value = 232088.4
# E12 intentionally left blank (matches E11, which is an empty string cell)
$ws.Range("F12").Value = 1085328.5
$ws.Range("G12").Value = 1639298.2
$ws.Range("H12").Value = 420444.2
$ws.Range("I12").Value = 245684.1
$ws.Range("J12").Value = 248385.5
$ws.Range("K12").Value = 592962.7
$ws.Range("L12").Value = 1451635.8
# M12 intentionally left blank (matches M11, which is an empty string cell)
$ws.Range("N12").Value = 28139
$ws.Range("O12").Value = 154145.6
$ws.Range("P12").Value = 365990.7
$ws.Range("Q12").Value = 2044598.5
$ws.Range("R12").Value = 3924.9
$ws.Range("S12").Value = 73526.5

$wb.Save()
